# Apply the "New simulation files for schemes report" edit.
#
# Summary of what changes on the single worksheet:
#  - Columns U:AD (which duplicated C:L under a second HKL block) are removed
#    entirely, shrinking the used range from A1:AD19 to A1:T29.
#  - Row 2 (the HKL/index header row) keeps its same layout but the [h k l]
#    labels in C2:T2 are reshuffled.
#  - The scheme-name labels in column B (rows 3-19) are replaced with a new
#    set of scheme names.
#  - Ten brand-new rows (20-29) are appended, continuing the same row layout
#    (index in A, scheme name in B, all 1's across C:T), with ten more new
#    scheme names.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. Drop the extra duplicated columns U:AD (rows 1-19) completely, so the
#    cells disappear from the sheet rather than merely going blank.
# ---------------------------------------------------------------------
$ws.Range("U1:AD19").Clear()

# ---------------------------------------------------------------------
# 2. Row 2 header labels (C2:T2) - same positions, new/reordered values.
# ---------------------------------------------------------------------
$row2Labels = @(
    "[3, 1, 0]",
    "[2, 2, 2]",
    "[1, 1, 0]",
    "[3, 2, 1]",
    "[4, 0, 0]",
    "[2, 1, 1]",
    "[2, 0, 0]",
    "[2, 2, 0]",
    "1Pair-A",
    "1Pair-B",
    "2Pairs-A",
    "2Pairs-B",
    "3Pairs-A",
    "3Pairs-B",
    "3Pairs-C",
    "4Pairs",
    "5A4F",
    "MaxUnique"
)
for ($i = 0; $i -lt $row2Labels.Length; $i++) {
    $ws.Cells.Item(2, 3 + $i).Value = $row2Labels[$i]
}

# ---------------------------------------------------------------------
# 3. Column B scheme names for the existing rows 3-19 get replaced with the
#    new scheme-name list (row/column layout and the all-1's data body are
#    unchanged).
# ---------------------------------------------------------------------
$existingRowNames = @(
    "Spiral5",
    "RotRing OmegaMax-90",
    "Equal Angle",
    "Tilt Rotate",
    "CLR",
    "Rizzie Hex",
    "Thomas Hex",
    "Tilt Rotate_Partial",
    "RotRing OmegaMax-60",
    "Equal Angle_Partial",
    "Rizzie Hex_Partial",
    "ND Single",
    "RD Single",
    "TD Single",
    "Morris Single",
    "Ring Perpendicular to ND",
    "Ring Perpendicular to RD"
)
for ($i = 0; $i -lt $existingRowNames.Length; $i++) {
    $ws.Cells.Item(3 + $i, 2).Value = $existingRowNames[$i]
}

# ---------------------------------------------------------------------
# 4. Append the ten brand-new rows 20-29, matching the formatting of the
#    rows directly above them (bold/bordered/centered column A, plain
#    column B and data cells).
# ---------------------------------------------------------------------
$newRowNames = @(
    "Ring Perpendicular to TD",
    "OffsetFTD",
    "OffsetATD",
    "OffsetF45",
    "OffsetA45",
    "OffsetFRD",
    "OffsetARD",
    "Gaussian Quadrature",
    "Michael-CCHex",
    "Michael-SNHex"
)

$ws.Range("A19:T19").Copy()
$ws.Range("A20:T29").PasteSpecial(-4122)  # xlPasteFormats

for ($i = 0; $i -lt $newRowNames.Length; $i++) {
    $r = 20 + $i
    $ws.Cells.Item($r, 1).Value = 18 + $i
    $ws.Cells.Item($r, 2).Value = $newRowNames[$i]
    for ($c = 3; $c -le 20; $c++) {
        $ws.Cells.Item($r, $c).Value = 1
    }
}

Write-Output "done"
